# Insert a new data row at row 198 (pushing the existing row 198..270 down
# to 199..271) and populate it with the new "Ají" price record for
# Feria Lagunitas de Puerto Montt, matching the weekly update commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 198:270 down by one row, leaving a blank row 198 behind
# (inherits formatting, e.g. the date style on column D, from the row above).
$ws.Rows.Item(198).Insert()

# Fill in the new record in the now-empty row 198.
$ws.Range("A198").Value = 4
$ws.Range("B198").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C198").Value = "Los Lagos"
$ws.Range("D198").Value = 44726
$ws.Range("E198").Value = 10
$ws.Range("F198").Value = 100112021
$ws.Range("G198").Value = "Ají"
$ws.Range("H198").Value = "Inferno"
$ws.Range("I198").Value = "Primera"
$ws.Range("J198").Value = 140
$ws.Range("K198").Value = 24000
$ws.Range("L198").Value = 25000
$ws.Range("M198").Value = 24500
$ws.Range("N198").Value = "`$/caja 12 kilos"
$ws.Range("O198").Value = "Región de Arica y Parinacota"
$ws.Range("P198").Value = 2042
$ws.Range("Q198").Value = 12
$ws.Range("R198").Value = "Hortaliza"
